# Add a new text box to slide 2 ("This is 1st page"), per the diff:
#   <p:sp> TextBox 2 at off(990600,3033962) ext(7391400,369332), noFill,
#   bodyPr wrap="square" with spAutoFit, text "Added this Text Box in
#   this slide." in red (FF0000).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# EMU -> points (1 pt = 12700 EMU): 990600/12700=78, 3033962/12700=238.8946...,
# 7391400/12700=582, 369332/12700=29.0812...
$shp = $s.Shapes.AddTextbox(1, 78, 238.89464566929135, 582, 29.081259842519685)
$shp.Name = "TextBox 2"
$shp.Fill.Visible = $false

$tf = $shp.TextFrame
$tf.WordWrap = $true
$tf.AutoSize = 1       # ppAutoSizeShapeToFitText -> <a:spAutoFit/>

$tr = $tf.TextRange
$tr.Text = "Added this Text Box in this slide."
$tr.Font.Color.RGB = 255   # RGB(255,0,0) red -> srgbClr FF0000
